{"js": "// Replace the date line and the twenty-five division expressions with\n// their updated values, per the commit's regenerated output.\nconst replacements = [\n  [\"2024-10-19 Saturday\", \"2024-10-20 Sunday\"],\n  [\"49\u00f79=\", \"55\u00f74=\"],\n  [\"50\u00f76=\", \"97\u00f75=\"],\n  [\"18\u00f77=\", \"82\u00f74=\"],\n  [\"88\u00f79=\", \"89\u00f73=\"],\n  [\"39\u00f76=\", \"62\u00f74=\"],\n  [\"27\u00f79=\", \"84\u00f76=\"],\n  [\"62\u00f78=\", \"22\u00f79=\"],\n  [\"59\u00f73=\", \"31\u00f79=\"],\n  [\"44\u00f77=\", \"11\u00f77=\"],\n  [\"43\u00f77=\", \"82\u00f72=\"],\n  [\"40\u00f76=\", \"54\u00f78=\"],\n  [\"22\u00f75=\", \"22\u00f73=\"],\n  [\"64\u00f74=\", \"86\u00f78=\"],\n  [\"68\u00f72=\", \"46\u00f72=\"],\n  [\"16\u00f75=\", \"57\u00f73=\"],\n  [\"66\u00f75=\", \"40\u00f77=\"],\n  [\"30\u00f77=\", \"71\u00f74=\"],\n  [\"83\u00f73=\", \"52\u00f75=\"],\n  [\"96\u00f79=\", \"13\u00f79=\"],\n  [\"47\u00f75=\", \"23\u00f77=\"],\n  [\"38\u00f75=\", \"98\u00f78=\"],\n  [\"93\u00f77=\", \"34\u00f79=\"],\n  [\"56\u00f73=\", \"71\u00f72=\"],\n  [\"67\u00f77=\", \"43\u00f73=\"],\n  [\"55\u00f75=\", \"74\u00f78=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and the twenty-five division expressions with\n# their updated values, per the commit's regenerated output.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-10-19 Saturday\", \"2024-10-20 Sunday\"),\n    @(\"49\u00f79=\", \"55\u00f74=\"),\n    @(\"50\u00f76=\", \"97\u00f75=\"),\n    @(\"18\u00f77=\", \"82\u00f74=\"),\n    @(\"88\u00f79=\", \"89\u00f73=\"),\n    @(\"39\u00f76=\", \"62\u00f74=\"),\n    @(\"27\u00f79=\", \"84\u00f76=\"),\n    @(\"62\u00f78=\", \"22\u00f79=\"),\n    @(\"59\u00f73=\", \"31\u00f79=\"),\n    @(\"44\u00f77=\", \"11\u00f77=\"),\n    @(\"43\u00f77=\", \"82\u00f72=\"),\n    @(\"40\u00f76=\", \"54\u00f78=\"),\n    @(\"22\u00f75=\", \"22\u00f73=\"),\n    @(\"64\u00f74=\", \"86\u00f78=\"),\n    @(\"68\u00f72=\", \"46\u00f72=\"),\n    @(\"16\u00f75=\", \"57\u00f73=\"),\n    @(\"66\u00f75=\", \"40\u00f77=\"),\n    @(\"30\u00f77=\", \"71\u00f74=\"),\n    @(\"83\u00f73=\", \"52\u00f75=\"),\n    @(\"96\u00f79=\", \"13\u00f79=\"),\n    @(\"47\u00f75=\", \"23\u00f77=\"),\n    @(\"38\u00f75=\", \"98\u00f78=\"),\n    @(\"93\u00f77=\", \"34\u00f79=\"),\n    @(\"56\u00f73=\", \"71\u00f72=\"),\n    @(\"67\u00f77=\", \"43\u00f73=\"),\n    @(\"55\u00f75=\", \"74\u00f78=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
